# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve
# profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2499.7778
$ws.Range("I40").Value = 2666.6667
$ws.Range("J40").Value = 2166
$ws.Range("K40").Value = 2666.6667
$ws.Range("L40").Value = 2166
$ws.Range("M40").Value = -2491.6667
$ws.Range("N40").Value = -2516
$ws.Range("H51").Value = 9500
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 13000
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 13000
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -13968
$ws.Range("H64").Value = 9499.182000000001
$ws.Range("I64").Value = 5562.125
$ws.Range("K64").Value = 5562.125
$ws.Range("M64").Value = -5314.125
$ws.Range("H67").Value = 9499.182000000001
$ws.Range("I67").Value = 5562.125
$ws.Range("K67").Value = 5562.125
$ws.Range("M67").Value = -4704.125
$ws.Range("H86").Value = 8932.5
$ws.Range("I86").Value = 8533.333000000001
$ws.Range("K86").Value = 8533.333000000001
$ws.Range("M86").Value = -7410.333000000001
$ws.Range("H88").Value = 2546.5625
$ws.Range("I88").Value = 2665
$ws.Range("J88").Value = 2519.2307
$ws.Range("K88").Value = 2665
$ws.Range("L88").Value = 2519.2307
$ws.Range("M88").Value = -2259
$ws.Range("N88").Value = -3331.2307
$ws.Range("H89").Value = 8932.5
$ws.Range("I89").Value = 8533.333000000001
$ws.Range("K89").Value = 42666.665
$ws.Range("M89").Value = -37050.665
$ws.Range("H91").Value = 2546.5625
$ws.Range("I91").Value = 2665
$ws.Range("J91").Value = 2519.2307
$ws.Range("K91").Value = 2665
$ws.Range("L91").Value = 2519.2307
$ws.Range("M91").Value = -1261
$ws.Range("N91").Value = -5327.2307
$ws.Range("H99").Value = 1626
$ws.Range("I99").Value = 1626
$ws.Range("K99").Value = 4878
$ws.Range("M99").Value = -3380
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H107").Value = 2326
$ws.Range("I107").Value = 2326
$ws.Range("K107").Value = 2326
$ws.Range("M107").Value = -406
$ws.Range("H116").Value = 5419.9
$ws.Range("I116").Value = 2750
$ws.Range("J116").Value = 6087.375
$ws.Range("K116").Value = 2750
$ws.Range("L116").Value = 6087.375
$ws.Range("M116").Value = 692
$ws.Range("N116").Value = -12971.375
$ws.Range("H127").Value = 2864.2
$ws.Range("I127").Value = 2864.2
$ws.Range("K127").Value = 8592.599999999999
$ws.Range("M127").Value = -3632.599999999999
$ws.Range("H129").Value = 3343.875
$ws.Range("I129").Value = 933.4286
$ws.Range("K129").Value = 2800.2858
$ws.Range("M129").Value = 2199.7142
$ws.Range("H132").Value = 2625.8076
$ws.Range("I132").Value = 2403.0952
$ws.Range("K132").Value = 7209.285600000001
$ws.Range("M132").Value = -4679.285600000001
$ws.Range("H138").Value = 2641.5
$ws.Range("I138").Value = 1681
$ws.Range("J138").Value = 4562.5
$ws.Range("K138").Value = 5043
$ws.Range("L138").Value = 13687.5
$ws.Range("M138").Value = 97
$ws.Range("N138").Value = -23967.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 5373.5
$ws.Range("I10").Value = 2247.5
$ws.Range("J10").Value = 8499.5
$ws.Range("K10").Value = 2247.5
$ws.Range("L10").Value = 8499.5
$ws.Range("M10").Value = -2077.5
$ws.Range("N10").Value = -8839.5
$ws.Range("H61").Value = 893.75
$ws.Range("I61").Value = 893.75
$ws.Range("K61").Value = 893.75
$ws.Range("M61").Value = -681.75
$ws.Range("H92").Value = 43999.5
$ws.Range("J92").Value = 43999.5
$ws.Range("L92").Value = 43999.5
$ws.Range("N92").Value = -48991.5
$ws.Range("H102").Value = 5142.7144
$ws.Range("I102").Value = 5142.7144
$ws.Range("K102").Value = 5142.7144
$ws.Range("M102").Value = -3520.7144
$ws.Range("H132").Value = 1999.5
$ws.Range("I132").Value = 1999.5
$ws.Range("K132").Value = 5998.5
$ws.Range("M132").Value = -3468.5
$ws.Range("H136").Value = 893.75
$ws.Range("I136").Value = 893.75
$ws.Range("K136").Value = 2681.25
$ws.Range("M136").Value = -131.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H29").Value = 1343.6666
$ws.Range("I29").Value = 1343.6666
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1343.6666
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1054.6666
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2109.8333
$ws.Range("I31").Value = 2074.4546
$ws.Range("K31").Value = 2074.4546
$ws.Range("M31").Value = -1779.4546
$ws.Range("H34").Value = 2109.8333
$ws.Range("I34").Value = 2074.4546
$ws.Range("K34").Value = 2074.4546
$ws.Range("M34").Value = -1872.4546
$ws.Range("H58").Value = 2088.4285
$ws.Range("I58").Value = 1523.8
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 1523.8
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -1320.8
$ws.Range("N58").Value = -3906
$ws.Range("H132").Value = 4559.8237
$ws.Range("I132").Value = 4534.7856
$ws.Range("K132").Value = 13604.3568
$ws.Range("M132").Value = -11074.3568
$ws.Range("H134").Value = 1605.6538
$ws.Range("I134").Value = 1225.75
$ws.Range("J134").Value = 2872
$ws.Range("K134").Value = 3677.25
$ws.Range("L134").Value = 8616
$ws.Range("M134").Value = -1142.25
$ws.Range("N134").Value = -13686
$ws.Range("H136").Value = 2088.4285
$ws.Range("I136").Value = 1523.8
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4571.4
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2021.4
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 5428.5713
$ws.Range("J49").Value = 5500
$ws.Range("L49").Value = 16500
$ws.Range("N49").Value = -16812
$ws.Range("H56").Value = 9997.666999999999
$ws.Range("I56").Value = 9997.666999999999
$ws.Range("K56").Value = 9997.666999999999
$ws.Range("M56").Value = -9467.666999999999
$ws.Range("H80").Value = 6832.6665
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 14998
$ws.Range("K80").Value = 8250
$ws.Range("L80").Value = 44994
$ws.Range("M80").Value = -7314
$ws.Range("N80").Value = -46866
$ws.Range("H83").Value = 6832.6665
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 14998
$ws.Range("K83").Value = 24750
$ws.Range("L83").Value = 134982
$ws.Range("M83").Value = -20070
$ws.Range("N83").Value = -144342
$ws.Range("H128").Value = 277759.5
$ws.Range("I128").Value = 277759.5
$ws.Range("K128").Value = 833278.5
$ws.Range("M128").Value = -828298.5
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 533333
$ws.Range("J24").Value = 38094.855
$ws.Range("L24").Value = 38094.855
$ws.Range("N24").Value = -38440.855
$ws.Range("H46").Value = 35780
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 35780
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 35780
$ws.Range("N46").Value = -36092
$ws.Range("M46").ClearContents()
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H97").Value = 1022
$ws.Range("I97").Value = 706.6667
$ws.Range("J97").Value = 1495
$ws.Range("K97").Value = 706.6667
$ws.Range("L97").Value = 1495
$ws.Range("M97").Value = -210.6667
$ws.Range("N97").Value = -2487

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 220.5
$ws.Range("I16").Value = 220.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 220.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -50.5
$ws.Range("N16").ClearContents()
$ws.Range("H62").Value = 47555
$ws.Range("J62").Value = 47555
$ws.Range("L62").Value = 47555
$ws.Range("N62").Value = -48927
$ws.Range("H65").Value = 47555
$ws.Range("J65").Value = 47555
$ws.Range("L65").Value = 142665
$ws.Range("N65").Value = -149529
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1500
$ws.Range("N82").Value = -2222
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1500
$ws.Range("N85").Value = -3996
$ws.Range("M85").ClearContents()
$ws.Range("H132").Value = 1879.091
$ws.Range("I132").Value = 1816.4445
$ws.Range("J132").Value = 2161
$ws.Range("K132").Value = 5449.333500000001
$ws.Range("L132").Value = 6483
$ws.Range("M132").Value = -2919.333500000001
$ws.Range("N132").Value = -11543

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -11058
$ws.Range("H22").Value = 17500
$ws.Range("J22").Value = 17500
$ws.Range("L22").Value = 17500
$ws.Range("N22").Value = -18086
$ws.Range("H132").Value = 1431.5454
$ws.Range("I132").Value = 1502.421
$ws.Range("K132").Value = 4507.263
$ws.Range("M132").Value = -1977.263
